$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3326.0476
$ws.Range("I15").Value = 3326.0476
$ws.Range("K15").Value = 9978.1428
$ws.Range("M15").Value = -9809.1428
$ws.Range("H18").Value = 769.3333
$ws.Range("I18").Value = 752.8570999999999
$ws.Range("K18").Value = 752.8570999999999
$ws.Range("M18").Value = -468.8570999999999
$ws.Range("H28").Value = 1428.0555
$ws.Range("I28").Value = 1756.1111
$ws.Range("J28").Value = 1100
$ws.Range("K28").Value = 1756.1111
$ws.Range("L28").Value = 1100
$ws.Range("M28").Value = -1271.1111
$ws.Range("N28").Value = -2070
$ws.Range("H43").Value = 1270.0625
$ws.Range("I43").Value = 499.75
$ws.Range("J43").Value = 1380.1072
$ws.Range("K43").Value = 499.75
$ws.Range("L43").Value = 1380.1072
$ws.Range("M43").Value = -430.75
$ws.Range("N43").Value = -1518.1072
$ws.Range("H86").Value = 1739.5
$ws.Range("I86").Value = 1581.8182
$ws.Range("J86").Value = 1987.2858
$ws.Range("K86").Value = 1581.8182
$ws.Range("L86").Value = 1987.2858
$ws.Range("M86").Value = -458.8181999999999
$ws.Range("N86").Value = -4233.2858
$ws.Range("H89").Value = 1739.5
$ws.Range("I89").Value = 1581.8182
$ws.Range("J89").Value = 1987.2858
$ws.Range("K89").Value = 7909.090999999999
$ws.Range("L89").Value = 9936.429
$ws.Range("M89").Value = -2293.090999999999
$ws.Range("N89").Value = -21168.429
$ws.Range("H92").Value = 3372.7368
$ws.Range("I92").Value = 3006.3333
$ws.Range("J92").Value = 4000.8572
$ws.Range("K92").Value = 3006.3333
$ws.Range("L92").Value = 4000.8572
$ws.Range("M92").Value = -1758.3333
$ws.Range("N92").Value = -6496.8572
$ws.Range("H98").Value = 1839.6957
$ws.Range("I98").Value = 2455.9092
$ws.Range("J98").Value = 1274.8334
$ws.Range("K98").Value = 2455.9092
$ws.Range("L98").Value = 1274.8334
$ws.Range("M98").Value = -957.9092000000001
$ws.Range("N98").Value = -4270.8334
$ws.Range("H99").Value = 1245.8889
$ws.Range("I99").Value = 475.16666
$ws.Range("K99").Value = 1425.49998
$ws.Range("M99").Value = 72.50001999999995
$ws.Range("H103").Value = 1150
$ws.Range("I103").Value = 725
$ws.Range("J103").Value = 1575
$ws.Range("K103").Value = 2175
$ws.Range("L103").Value = 4725
$ws.Range("M103").Value = -1589
$ws.Range("N103").Value = -5897
$ws.Range("H106").Value = 3876.1
$ws.Range("I106").Value = 4234
$ws.Range("J106").Value = 2802.4
$ws.Range("K106").Value = 4234
$ws.Range("L106").Value = 2802.4
$ws.Range("M106").Value = -3603
$ws.Range("N106").Value = -4064.4
$ws.Range("H112").Value = 23257430
$ws.Range("J112").Value = 1818.0541
$ws.Range("L112").Value = 5454.1623
$ws.Range("N112").Value = -7670.1623
$ws.Range("H122").Value = 1839.6957
$ws.Range("I122").Value = 2455.9092
$ws.Range("J122").Value = 1274.8334
$ws.Range("K122").Value = 7367.7276
$ws.Range("L122").Value = 3824.5002
$ws.Range("M122").Value = -4917.7276
$ws.Range("N122").Value = -8724.5002

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4080275.8
$ws.Range("I32").Value = 4645767
$ws.Range("K32").Value = 4645767
$ws.Range("M32").Value = -4645480
$ws.Range("H76").Value = 39972.8
$ws.Range("J76").Value = 39972.8
$ws.Range("L76").Value = 39972.8
$ws.Range("N76").Value = -40648.8
$ws.Range("H79").Value = 39972.8
$ws.Range("J79").Value = 39972.8
$ws.Range("L79").Value = 39972.8
$ws.Range("N79").Value = -42312.8
$ws.Range("H106").Value = 39391.43
$ws.Range("J106").Value = 39391.43
$ws.Range("L106").Value = 39391.43
$ws.Range("N106").Value = -41915.43
$ws.Range("H122").Value = 6946913
$ws.Range("I122").Value = 2753.4614
$ws.Range("J122").Value = 37038270
$ws.Range("K122").Value = 8260.3842
$ws.Range("L122").Value = 111114810
$ws.Range("M122").Value = -5810.3842
$ws.Range("N122").Value = -111119710
$ws.Range("H132").Value = 47621.09
$ws.Range("I132").Value = 38976.258
$ws.Range("J132").Value = 61351.117
$ws.Range("K132").Value = 116928.774
$ws.Range("L132").Value = 184053.351
$ws.Range("M132").Value = -114398.774
$ws.Range("N132").Value = -189113.351

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 466.77777
$ws.Range("I37").Value = 275.125
$ws.Range("K37").Value = 275.125
$ws.Range("M37").Value = -138.125
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 33462.812
$ws.Range("I132").Value = 1911.3334
$ws.Range("J132").Value = 74029
$ws.Range("K132").Value = 5734.0002
$ws.Range("L132").Value = 222087
$ws.Range("M132").Value = -3204.0002
$ws.Range("N132").Value = -227147

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 816.9663
$ws.Range("I68").Value = 446.64615
$ws.Range("J68").Value = 1819.9166
$ws.Range("K68").Value = 1339.93845
$ws.Range("L68").Value = 5459.7498
$ws.Range("M68").Value = -528.9384499999999
$ws.Range("N68").Value = -7081.7498
$ws.Range("H71").Value = 816.9663
$ws.Range("I71").Value = 446.64615
$ws.Range("J71").Value = 1819.9166
$ws.Range("K71").Value = 4019.81535
$ws.Range("L71").Value = 16379.2494
$ws.Range("M71").Value = 36.18465000000015
$ws.Range("N71").Value = -24491.2494
$ws.Range("H107").Value = 926.62195
$ws.Range("I107").Value = 373.53333
$ws.Range("J107").Value = 1599.2972
$ws.Range("K107").Value = 1120.59999
$ws.Range("L107").Value = 4797.8916
$ws.Range("M107").Value = 799.4000100000001
$ws.Range("N107").Value = -8637.891599999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 6386.8
$ws.Range("I43").Value = 644.6667
$ws.Range("J43").Value = 15000
$ws.Range("K43").Value = 644.6667
$ws.Range("L43").Value = 15000
$ws.Range("M43").Value = -493.6667
$ws.Range("N43").Value = -15302
$ws.Range("H46").Value = 26000
$ws.Range("J46").Value = 26000
$ws.Range("L46").Value = 26000
$ws.Range("N46").Value = -26312
$ws.Range("H57").Value = 15000
$ws.Range("J57").Value = 15000
$ws.Range("L57").Value = 15000
$ws.Range("N57").Value = -16640
$ws.Range("H102").Value = 2018.3077
$ws.Range("I102").Value = 2002.4
$ws.Range("J102").Value = 2071.3333
$ws.Range("K102").Value = 2002.4
$ws.Range("L102").Value = 2071.3333
$ws.Range("M102").Value = -380.4000000000001
$ws.Range("N102").Value = -5315.3333
$ws.Range("H122").Value = 1907.1428
$ws.Range("I122").Value = 1391.6666
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 4174.9998
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -1724.9998
$ws.Range("N122").Value = -19900
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 882.8946999999999
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 1106.25
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 1106.25
$ws.Range("M22").Value = -205
$ws.Range("N22").Value = -1696.25
$ws.Range("H27").Value = 882.8946999999999
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 1106.25
$ws.Range("K27").Value = 500
$ws.Range("L27").Value = 1106.25
$ws.Range("M27").Value = -393
$ws.Range("N27").Value = -1320.25
$ws.Range("H46").Value = 542.8570999999999
$ws.Range("I46").Value = 520
$ws.Range("K46").Value = 520
$ws.Range("M46").Value = -332
$ws.Range("H55").Value = 272.9643
$ws.Range("I55").Value = 233.23077
$ws.Range("J55").Value = 307.4
$ws.Range("K55").Value = 233.23077
$ws.Range("L55").Value = 307.4
$ws.Range("M55").Value = -60.23077000000001
$ws.Range("N55").Value = -653.4
$ws.Range("H68").Value = 1745.4546
$ws.Range("I68").Value = 1657.1428
$ws.Range("J68").Value = 1900
$ws.Range("K68").Value = 1657.1428
$ws.Range("L68").Value = 1900
$ws.Range("M68").Value = -908.1428000000001
$ws.Range("N68").Value = -3398
$ws.Range("H71").Value = 1745.4546
$ws.Range("I71").Value = 1657.1428
$ws.Range("J71").Value = 1900
$ws.Range("K71").Value = 8285.714
$ws.Range("L71").Value = 9500
$ws.Range("M71").Value = -4541.714
$ws.Range("N71").Value = -16988
$ws.Range("H76").Value = 35500
$ws.Range("J76").Value = 35500
$ws.Range("L76").Value = 35500
$ws.Range("N76").Value = -36176
$ws.Range("H79").Value = 35500
$ws.Range("J79").Value = 35500
$ws.Range("L79").Value = 35500
$ws.Range("N79").Value = -37840
$ws.Range("H122").Value = 3159.8
$ws.Range("I122").Value = 2950
$ws.Range("J122").Value = 3999
$ws.Range("K122").Value = 8850
$ws.Range("L122").Value = 11997
$ws.Range("M122").Value = -6400
$ws.Range("N122").Value = -16897

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H127").Value = 50429
$ws.Range("J127").Value = 50429
$ws.Range("L127").Value = 50429
$ws.Range("N127").Value = -60349
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H132").Value = 79377.42
$ws.Range("I132").Value = 65019.75
$ws.Range("J132").Value = 102349.7
$ws.Range("K132").Value = 195059.25
$ws.Range("L132").Value = 307049.1
$ws.Range("M132").Value = -192529.25
$ws.Range("N132").Value = -312109.1

Write-Output "Applied all cell updates."